# Large commit of various bits:
#  - Fill in a previously-blank header row (row 27) with column labels.
#  - Update the saved window view: scroll position and active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 was an empty row between the regression-coefficient table (ending
# row 26) and the "by group/year" table (starting row 28). Add header
# labels for that second table.
$ws.Range("A27").Value2 = "Group"
$ws.Range("B27").Value2 = "Year"
$ws.Range("C27").Value2 = "Coefficient"
$ws.Range("D27").Value2 = "Standard error"
$ws.Range("E27").Value2 = "T value"

# Update the active window's scroll position and selection to match the
# new view of the sheet (scrolled up a bit, with E28 now selected).
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("E28").Select()
